# Keep section and recipes on same page in table of contents.
# Insert a new ingredient row ("fennel") at row 68 of the Ingredients sheet,
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ingredients")

# Insert a new row above the current row 68 ("fish fingers"), pushing it
# (and everything below it) down by one row.
$ws.Rows.Item(68).Insert()

# Populate the new row with the fennel ingredient data.
$ws.Range("A68").Value = "fennel"
$ws.Range("B68").Value = "Vegetable"
$ws.Range("C68").Value = 14
$ws.Range("D68").Value = 1
$ws.Range("E68").Value = 2
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0

# Keep the selection/active cell on the newly inserted row.
$ws.Range("G68").Select()

# The filter database named range needs to grow by one row to keep covering
# the full table (it previously ended at row 178, now 179).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Ingredients!_FilterDatabase") {
        $n.RefersTo = "=Ingredients!`$A`$2:`$G`$179"
    }
}
